# Generate Report for Archive
# Update the "Ready for handoff" status text (wherever it appears) to
# "In Translation", then re-autofit the columns that hold that status
# text so the column widths reflect the new (shorter) label.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Columns.Item(5).AutoFit() | Out-Null
$ws1.Columns.Item(6).AutoFit() | Out-Null

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Columns.Item(3).AutoFit() | Out-Null

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Columns.Item(3).AutoFit() | Out-Null
